# Updated capital structure database
# Applies the refreshed australia_insurance_general figures to rows 2-6
# (row order for AUB / PSC / Steadfast also rotates; see company_name updates below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - industry aggregate (col B = "4")
$ws.Range("D2").Value = 0.1513
$ws.Range("E2").Value = 0.1453
$ws.Range("G2").Value = 0.2841271854717659
$ws.Range("H2").Value = 0.2841271854717659
$ws.Range("I2").Value = 0.246083129956427
$ws.Range("J2").Value = 0.189021862846119
$ws.Range("K2").Value = 5.31
$ws.Range("L2").Value = 0.005612989154563329
$ws.Range("M2").Value = 79.99000000000001
$ws.Range("N2").Value = 0.01850699304744167
$ws.Range("O2").Value = 15.06403013182674
$ws.Range("P2").Value = 76.5
$ws.Range("Q2").Value = 0.01769952454218386
$ws.Range("R2").Value = 14.40677966101695
$ws.Range("S2").Value = 3.490000000000002
$ws.Range("T2").Value = 0.04363045380672586
$ws.Range("U2").Value = 222.481
$ws.Range("V2").Value = 0.05147461332901448
$ws.Range("W2").Value = 0.1149025495840194
$ws.Range("X2").Value = 0.04481290168218622
$ws.Range("Y2").Value = 0.07008964790183317
$ws.Range("Z2").Value = 2.696272311665755
$ws.Range("AA2").Value = 0.3615624893741434
$ws.Range("AB2").Value = 0.03969251258579146
$ws.Range("AC2").Value = 0.3218699767883519
$ws.Range("AD2").Value = 830.91
$ws.Range("AE2").Value = 3.152186993104562
$ws.Range("AF2").Value = 834.0621869931045
$ws.Range("AG2").Value = 611.5811869931046
$ws.Range("AH2").Value = 0.1617587012995863
$ws.Range("AI2").Value = 0.3846392408267462
$ws.Range("AJ2").Value = 0.1239591627134912
$ws.Range("AK2").Value = 0.314284737718326
$ws.Range("AL2").Value = 20.552
$ws.Range("AM2").Value = 19.472
$ws.Range("AN2").Value = 3.075621853716316
$ws.Range("AO2").Value = 11.29135850525496
$ws.Range("AP2").Value = 2.263774011671249
$ws.Range("AQ2").Value = 11.91762530813476

# Row 3 - AUB Group Limited (ASX:AUB)
$ws.Range("B3").Value = 'AUB Group Limited (ASX:AUB)'
$ws.Range("D3").Value = 0.0906
$ws.Range("E3").Value = 0.0626
$ws.Range("G3").Value = 0.2686825053995681
$ws.Range("H3").Value = 0.2686825053995681
$ws.Range("I3").Value = 0.2399117175005576
$ws.Range("J3").Value = 0.1996684616617544
$ws.Range("K3").Value = 32.6
$ws.Range("L3").Value = 0.1408207343412527
$ws.Range("M3").Value = 14.2
$ws.Range("N3").Value = 0.01528525296017223
$ws.Range("O3").Value = 0.4355828220858896
$ws.Range("P3").Value = 14.2
$ws.Range("Q3").Value = 0.01528525296017223
$ws.Range("R3").Value = 0.4355828220858896
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 58.3
$ws.Range("V3").Value = 0.0627556512378902
$ws.Range("W3").Value = 0.1119890072140158
$ws.Range("X3").Value = 0.04481332692775535
$ws.Range("Y3").Value = 0.06717568028626045
$ws.Range("Z3").Value = 2.91370204850486
$ws.Range("AA3").Value = 0.5817744057656679
$ws.Range("AB3").Value = 0.03963490713863608
$ws.Range("AC3").Value = 0.5421394986270318
$ws.Range("AD3").Value = 180.5
$ws.Range("AE3").Value = 3.152186993104562
$ws.Range("AF3").Value = 183.6521869931046
$ws.Range("AG3").Value = 125.3521869931046
$ws.Range("AH3").Value = 0.1650580380284128
$ws.Range("AI3").Value = 0.3509821679837999
$ws.Range("AJ3").Value = 0.1188902423113429
$ws.Range("AK3").Value = 0.2696023171839894
$ws.Range("AL3").Value = 5.89
$ws.Range("AM3").Value = 5.89
$ws.Range("AN3").Value = 2.8393896492056
$ws.Range("AO3").Value = 9.303904923599321
$ws.Range("AP3").Value = 1.971876466778426
$ws.Range("AQ3").Value = 9.303904923599321

# Row 4 - PSC Insurance Group Limited (ASX:PSI)
$ws.Range("B4").Value = 'PSC Insurance Group Limited (ASX:PSI)'
$ws.Range("D4").Value = 0.287
$ws.Range("E4").Value = 0.228
$ws.Range("G4").Value = 0.3094841930116473
$ws.Range("H4").Value = 0.3094841930116473
$ws.Range("I4").Value = 0.2612312811980033
$ws.Range("J4").Value = 0.193371201254614
$ws.Range("K4").Value = 12.3
$ws.Range("L4").Value = 0.1023294509151414
$ws.Range("M4").Value = 15.3
$ws.Range("N4").Value = 0.02115304852758192
$ws.Range("O4").Value = 1.24390243902439
$ws.Range("P4").Value = 15.3
$ws.Range("Q4").Value = 0.02115304852758192
$ws.Range("R4").Value = 1.24390243902439
$ws.Range("U4").Value = 17.9
$ws.Range("V4").Value = 0.0247476842250795
$ws.Range("W4").Value = 0.117816091954023
$ws.Range("X4").Value = 0.04414042862424195
$ws.Range("Y4").Value = 0.07367566332978104
$ws.Range("Z4").Value = 1.917065390749601
$ws.Range("AA4").Value = 0.3707052374928962
$ws.Range("AB4").Value = 0.03974997633165189
$ws.Range("AC4").Value = 0.3309552611612444
$ws.Range("AD4").Value = 120.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 120.7
$ws.Range("AG4").Value = 102.8
$ws.Range("AH4").Value = 0.1430094786729858
$ws.Range("AI4").Value = 0.4169257340241797
$ws.Range("AJ4").Value = 0.1244401404188355
$ws.Range("AK4").Value = 0.3784977908689249
$ws.Range("AL4").Value = 4.7
$ws.Range("AM4").Value = 3.62
$ws.Range("AN4").Value = 3.244623655913978
$ws.Range("AO4").Value = 6.680851063829786
$ws.Range("AP4").Value = 2.763440860215054
$ws.Range("AQ4").Value = 8.674033149171271

# Row 5 - Steadfast Group Limited (ASX:SDF)
$ws.Range("B5").Value = 'Steadfast Group Limited (ASX:SDF)'
$ws.Range("D5").Value = 0.212
$ws.Range("G5").Value = 0.2883217846881865
$ws.Range("H5").Value = 0.2883217846881865
$ws.Range("I5").Value = 0.2486057123542335
$ws.Range("J5").Value = 0.1243028561771168
$ws.Range("K5").Value = -38.1
$ws.Range("L5").Value = -0.06439073854994085
$ws.Range("M5").Value = 50.49
$ws.Range("N5").Value = 0.0189612438035151
$ws.Range("O5").Value = -1.325196850393701
$ws.Range("P5").Value = 47
$ws.Range("Q5").Value = 0.01765059336037254
$ws.Range("R5").Value = -1.233595800524934
$ws.Range("S5").Value = 3.490000000000002
$ws.Range("T5").Value = 0.06912259853436328
$ws.Range("U5").Value = 145.4
$ws.Range("V5").Value = 0.05460417605528015
$ws.Range("W5").Value = -0.05351123595505618
$ws.Range("X5").Value = 0.04481247643661709
$ws.Range("Y5").Value = -0.09832371239167328
$ws.Range("Z5").Value = 2.835170100622905
$ws.Range("AA5").Value = 0.3524197412553906
$ws.Range("AB5").Value = 0.03963504883993102
$ws.Range("AC5").Value = 0.3127846924154596
$ws.Range("AD5").Value = 526.3
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 526.3
$ws.Range("AG5").Value = 380.9
$ws.Range("AH5").Value = 0.1650308864569941
$ws.Range("AI5").Value = 0.3889874353288987
$ws.Range("AJ5").Value = 0.1251437395275487
$ws.Range("AK5").Value = 0.3154190129181849
$ws.Range("AL5").Value = 9.449999999999999
$ws.Range("AM5").Value = 9.449999999999999
$ws.Range("AN5").Value = 3.084994138335287
$ws.Range("AO5").Value = 15.56613756613757
$ws.Range("AP5").Value = 2.232708089097303
$ws.Range("AQ5").Value = 15.56613756613757
$ws.Range("E5").ClearContents()

# Row 6 - Ensurance Limited (ASX:ENA)
$ws.Range("D6").Value = 0.07200000000000001
$ws.Range("G6").Value = -0.4618320610687023
$ws.Range("H6").Value = -0.4618320610687023
$ws.Range("I6").Value = -0.4732824427480916
$ws.Range("J6").Value = -0.4732824427480916
$ws.Range("K6").Value = -1.49
$ws.Range("L6").Value = -0.5687022900763359
$ws.Range("U6").Value = 0.881
$ws.Range("V6").Value = 0.1249645390070922
$ws.Range("W6").Value = 0.9254658385093167
$ws.Range("X6").Value = 0.05105882978849452
$ws.Range("Y6").Value = 0.8744070087208222
$ws.Range("Z6").Value = 262.0000000000056
$ws.Range("AA6").Value = -124.0000000000026
$ws.Range("AB6").Value = 0.04125497227618417
$ws.Range("AC6").Value = -124.0412549722788
$ws.Range("AD6").Value = 3.41
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 3.41
$ws.Range("AG6").Value = 2.529
$ws.Range("AH6").Value = 0.3260038240917782
$ws.Range("AI6").Value = 1.274766355140187
$ws.Range("AJ6").Value = 0.2640150328844347
$ws.Range("AK6").Value = 1.409698996655518
$ws.Range("AL6").Value = 0.512
$ws.Range("AM6").Value = 0.512
$ws.Range("AN6").Value = -2.818181818181818
$ws.Range("AO6").Value = -2.421875
$ws.Range("AP6").Value = -2.090082644628099
$ws.Range("AQ6").Value = -2.421875
